$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Game_Record")

# Row 16
$ws.Cells.Item(16, 1).Formula = "=ROW()-1"
$ws.Cells.Item(16, 2).Value = (Get-Date -Year 2025 -Month 10 -Day 26).Date
$ws.Cells.Item(16, 3).Value = "SiderFace"
$ws.Cells.Item(16, 4).Value = "DrSystomatix"
$ws.Cells.Item(16, 5).Value = "SimpleJack"
$ws.Cells.Item(16, 6).Value = "Player1"

# Row 17
$ws.Cells.Item(17, 1).Formula = "=ROW()-1"
$ws.Cells.Item(17, 2).Value = (Get-Date -Year 2025 -Month 11 -Day 21).Date
$ws.Cells.Item(17, 3).Value = "SiderFace"
$ws.Cells.Item(17, 4).Value = "Player1"
$ws.Cells.Item(17, 5).Value = "Doanage"
$ws.Cells.Item(17, 6).Value = "SimpleJack"

$ws.Range("G17").Select()
